# "Committed Corporate Customer excel file"
#
# The Create_Term_Deposits_LCY sheet already has 5 columns of header data
# (CUSTOMER.ID, CATEGORY, FTD.TYPE, PRINCIPAL, AUTO.ROLLOVER in A1:E1, with
# one data row in A2:E2). This adds the remaining Term-Deposit field names
# used by the Corporate Customer upload as new header columns F1:O1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "PROFIT.PAY.TERM",
    "INTEND.DATE",
    "CUST.REMARKS:1",
    "TAX.INTEREST.TYPE:1",
    "DRAWDOWN.ACCOUNT",
    "PRIN.LIQ.ACCT",
    "INT.LIQ.ACCT",
    "CHRG.LIQ.ACCT",
    "FINAL.MATURITY",
    "EXP.DATE"
)

$col = 6
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# Size the new columns to fit their header text, like the existing A:E
# columns (which are all best-fit/custom width).
$ws.Columns.Item(6).ColumnWidth = 15.917
$ws.Columns.Item(7).ColumnWidth = 11.917
$ws.Columns.Item(8).ColumnWidth = 15.251
$ws.Columns.Item(9).ColumnWidth = 18.75
$ws.Columns.Item(10).ColumnWidth = 21.417
$ws.Columns.Item(11).ColumnWidth = 13.084
$ws.Columns.Item(12).ColumnWidth = 11.584
$ws.Columns.Item(13).ColumnWidth = 13.584
$ws.Columns.Item(14).ColumnWidth = 14.917
$ws.Columns.Item(15).ColumnWidth = 8.584

# Leave the selection where the author left it when they were done (K3).
$ws.Range("K3").Select()
